# coeff react by coolant density
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new row at 9 (pushes old rows 9-15 down to 10-16).
#    This turns the old "header@10 / subheader@11 / data@12-15" block into
#    "header@11 / subheader@12 / data@13-16", matching the target layout,
#    and frees up row 9 for the new average-formula row.
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# Extend the coolant-temperature block header merge from B11:G11 to
# B11:H11 (it now needs to cover the extra alphaTtep column, H).
$ws.Range("B11:G11").UnMerge()
$ws.Range("B11:H11").Merge()

# ---------------------------------------------------------------------------
# 2. New shared strings must be introduced in the same order they first
#    appear in the target file, so that they land at sharedStrings indices
#    11, 12, 13 (matching the diff): the new section title, "alphaRotep",
#    then "react" (which replaces "ro" as a column caption in three spots).
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "Расчет коэффициента реактивности по плотности теплоносителя(Getera)"
$ws.Range("G20").Value = "alphaRotep"
$ws.Range("D4").Value = "react"
$ws.Range("L4").Value = "react"
$ws.Range("G12").Value = "react"
$ws.Range("L12").Value = "react"
$ws.Range("F20").Value = "react"
$ws.Range("L20").Value = "react"

# ---------------------------------------------------------------------------
# 3. Row 9: average of the fuel-temperature "alpha" columns (E and M)
# ---------------------------------------------------------------------------
$ws.Range("E9").Formula = "=SUM(E6:E8) / 3"
$ws.Range("M9").Formula = "=SUM(M6:M8) / 3"
$ws.Range("E9").Interior.Color = 5296274
$ws.Range("M9").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 4. Row 17: average of the coolant-temperature "alpha" columns (H and M)
# ---------------------------------------------------------------------------
$ws.Range("H17").Formula = "=SUM(H14:H16) / 3"
$ws.Range("M17").Formula = "=SUM(M14:M16) / 3"
$ws.Range("H17").Interior.Color = 5296274
$ws.Range("M17").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 5. New section (rows 19-25): reactivity coefficient by coolant density.
# ---------------------------------------------------------------------------

# Section title row (row 19), copied formatting from the row-11 title band.
$ws.Range("B11:H11").Copy()
$ws.Range("B19:H19").PasteSpecial(-4122)
$ws.Range("B19:G19").Merge()
$ws.Range("J19:P19").Merge()
$ws.Range("J19").Value = "Расчет коэффициента реактиввности по температуре теплоносителя(Scetch)"

# Sub-header row (row 20)
$ws.Range("B20").Value = "ro"
$ws.Range("C20").Value = "N(H)"
$ws.Range("D20").Value = "N(O)"
$ws.Range("E20").Value = "keff"
$ws.Range("J20").Value = "ro"
$ws.Range("K20").Value = "keff"
$ws.Range("M20").Value = "alphaTtep"
$ws.Range("G13").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("M12").Copy()
$ws.Range("M20").PasteSpecial(-4122)

# Data rows 21-24
$ws.Range("B21").Value = 0.69799999999999995
$ws.Range("C21").Formula = "=B21*6.02*10^23 / 18 *2 *1E-24"
$ws.Range("D21").Formula = "=B21*6.02*10^23 / 18 * 1E-24"
$ws.Range("E21").Value = 1.11696
$ws.Range("F21").Formula = "=(E21 - 1) / E21"
$ws.Range("J21").Value = 0.69799999999999995
$ws.Range("K21").Value = 1.018429
$ws.Range("L21").Formula = "=(K21 - 1) / K21"

$ws.Range("B22").Formula = "=1.1*B21"
$ws.Range("C22").Formula = "=B22*6.02*10^23 / 18 *2 *1E-24"
$ws.Range("D22").Formula = "=B22*6.02*10^23 / 18 * 1E-24"
$ws.Range("E22").Value = 1.1327100000000001
$ws.Range("F22").Formula = "=(E22 - 1) / E22"
$ws.Range("G22").Formula = "=(F21 - F22) / (B21 - B22)"
$ws.Range("J22").Formula = "=1.1*J21"
$ws.Range("K22").Value = 1.0226740000000001
$ws.Range("L22").Formula = "=(K22 - 1) / K22"
$ws.Range("M22").Formula = "=(L21-L22)/(J21-J22)"

$ws.Range("B23").Formula = "=1.2*B21"
$ws.Range("C23").Formula = "=B23*6.02*10^23 / 18 *2 *1E-24"
$ws.Range("D23").Formula = "=B23*6.02*10^23 / 18 * 1E-24"
$ws.Range("E23").Value = 1.12093
$ws.Range("F23").Formula = "=(E23 - 1) / E23"
$ws.Range("G23").Formula = "=(F22 - F23) / (B22 - B23)"
$ws.Range("J23").Formula = "=1.2*J21"
$ws.Range("K23").Value = 1.026467
$ws.Range("L23").Formula = "=(K23 - 1) / K23"
$ws.Range("M23").Formula = "=(L22-L23)/(J22-J23)"

$ws.Range("B24").Formula = "=1.3*B21"
$ws.Range("C24").Formula = "=B24*6.02*10^23 / 18 *2 *1E-24"
$ws.Range("D24").Formula = "=B24*6.02*10^23 / 18 * 1E-24"
$ws.Range("E24").Value = 1.1236900000000001
$ws.Range("F24").Formula = "=(E24 - 1) / E24"
$ws.Range("G24").Formula = "=(F23 - F24) / (B23 - B24)"
$ws.Range("J24").Formula = "=1.3*J21"
$ws.Range("K24").Value = 1.0298940000000001
$ws.Range("L24").Formula = "=(K24 - 1) / K24"
$ws.Range("M24").Formula = "=(L23-L24)/(J23-J24)"

# Formats for rows 21-24 (copy from the analogous rows 13-16 block)
$ws.Range("M13").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("M21").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("G22:G24").PasteSpecial(-4122)
$ws.Range("M14").Copy()
$ws.Range("M22:M24").PasteSpecial(-4122)
$ws.Range("G21").ClearContents()

# Row 25: average of the new density "alpha" columns (G and M)
$ws.Range("G25").Formula = "=SUM(G22:G24) / 3"
$ws.Range("M25").Formula = "=SUM(M22:M24) / 3"
$ws.Range("G25").Interior.Color = 5296274
$ws.Range("M25").Interior.Color = 5296274

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Selection / view bookkeeping (cosmetic, matches the authored commit).
# ---------------------------------------------------------------------------
$ws.Range("J26").Select()
